$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Merge the three "Poisson imag" / "e" / " editing" hyperlink runs
#    into a single run "Poisson image editing" while preserving the
#    hyperlink's run formatting (blue colour + single underline).
# ------------------------------------------------------------------
$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Replacement.ClearFormatting()
$rng.Find.Replacement.Font.Color = 13391121   # RGB(0x11,0x55,0xCC) -> BGR int
$rng.Find.Replacement.Font.Underline = 1      # wdUnderlineSingle
$rng.Find.Execute("Poisson image editing", $true, $false, $false, $false, $false, $true, 1, $false, "Poisson image editing", 2) | Out-Null

# ------------------------------------------------------------------
# 2) Highlight four phrases inside the Task 1 paragraph in yellow.
# ------------------------------------------------------------------
$phrases = @("smooth regions", "with edges", "the size", "increases")
foreach ($phrase in $phrases) {
    $rng2 = $d.Content
    $rng2.Find.ClearFormatting()
    $rng2.Find.Replacement.ClearFormatting()
    $rng2.Find.Replacement.Highlight = $true
    $rng2.Find.Execute($phrase, $true, $false, $false, $false, $false, $true, 1, $false, $phrase, 2) | Out-Null
}
